# "Add files via upload" — extends Sheet1 with three new "Case" columns
# (Case 16, Case 17, Case 18) that sit right after the existing Case 15
# column (P), mirroring the existing layout: a text header in row 1 and
# numeric values in rows 2-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header labels for the new columns Q, R, S (row 1)
$ws.Range("Q1").Value = "Case 16"
$ws.Range("R1").Value = "Case 17"
$ws.Range("S1").Value = "Case 18"

# Data values for rows 2-13 across the new columns Q, R, S
$data = @(
    @(671, 393, 255),
    @(1222, 331, 498),
    @(621, 264, 324),
    @(674, 327, 389),
    @(483, 324, 247),
    @(533, 171, 275),
    @(687, 171, 343),
    @(907, 223, 404),
    @(546, 222, 373),
    @(435, 219, 393),
    @(327, 203, 289),
    @(436, 214, 370)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 17).Value = $vals[0]
    $ws.Cells.Item($row, 18).Value = $vals[1]
    $ws.Cells.Item($row, 19).Value = $vals[2]
}

# Match the look of the existing data grid: copy the border formatting
# from the last existing column (P, "Case 15") onto the new Q:S block.
$ws.Range("P1:P13").Copy()
$ws.Range("Q1:S13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to match where editing left off.
$ws.Range("Q19").Select()
